$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.555.65"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.628.32"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  -0.18%  "

Set-TextValue $ws.Range("D5") "212.23"
$ws.Range("E5").Value = "  -0.10%  "

Set-TextValue $ws.Range("D6") "0.521"
$ws.Range("E6").Value = "  -0.23%  "

Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.22%  "

Set-TextValue $ws.Range("D8") "23.35"
$ws.Range("E8").Value = "  +1.83%  "

Set-TextValue $ws.Range("D9") "0.264"
$ws.Range("E9").Value = "  +2.50%  "

Set-TextValue $ws.Range("D10") "0.0612"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("D12").Value = "1.857.47"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("D13").Value = "1.629.12"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("E15").Value = "  -1.24%  "

Set-TextValue $ws.Range("D16") "65.40"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").Value = "27.523.02"
$ws.Range("E17").Value = "  -0.30%  "

Set-TextValue $ws.Range("D18") "230.27"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("E20").Value = "  -2.35%  "

Set-TextValue $ws.Range("D22") "10.53"
$ws.Range("E22").Value = "  +5.30%  "

Set-TextValue $ws.Range("D23") "4.36"
$ws.Range("E23").Value = "  +1.84%  "

Set-TextValue $ws.Range("D24") "2.15"
$ws.Range("E24").Value = "  +9.03%  "

Set-TextValue $ws.Range("D25") "149.19"
$ws.Range("E25").Value = "  -0.72%  "

Set-TextValue $ws.Range("D26") "6.89"
$ws.Range("E26").Value = "  -0.42%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").Value = "1.467.40"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("E34").Value = "  -1.82%  "

Set-TextValue $ws.Range("D35") "1.56"
$ws.Range("E35").Value = "  -1.00%  "

Set-TextValue $ws.Range("D36") "2.34"
$ws.Range("E36").Value = "  -1.80%  "

$ws.Range("E37").Value = "  +5.81%  "

Set-TextValue $ws.Range("D38") "0.879"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("E39").Value = "  +0.57%  "

Set-TextValue $ws.Range("D40") "0.554"
$ws.Range("E40").Value = "  -1.52%  "

$ws.Range("E41").Value = "  +1.93%  "

Set-TextValue $ws.Range("D42") "1.00"
$ws.Range("E42").Value = "  -0.14%  "

Set-TextValue $ws.Range("D43") "67.85"
$ws.Range("E43").Value = "  -2.81%  "

$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("E46").Value = "  -4.19%  "

Set-TextValue $ws.Range("D47") "1.76"
$ws.Range("E47").Value = "  +3.36%  "

$ws.Range("D48").Value = "1.767.16"
$ws.Range("E48").Value = "  -0.63%  "

Set-TextValue $ws.Range("D49") "87.45"
$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("E51").Value = "  +1.08%  "

# Row 28 and 29: coin data swap (BinanceUSD <-> EthereumClassic)
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "15.54"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  -0.14%  "